$wb = $excel.ActiveWorkbook

# Existing sheets: [1]="总计"  [2]="2022-Q3" (sheetId=2, holds old Q3 fund table)
$wsTotal = $wb.Worksheets.Item(1)
$wsQ3 = $wb.Worksheets.Item(2)

# 1) Duplicate the "2022-Q3" sheet right after itself. The duplicate keeps the
#    original fund-table contents/styles untouched and becomes the new
#    "2022-Q3" tab (new sheetId), while the original sheet (sheetId kept) is
#    turned into the new "2022-Q4" tab below.
$wsQ3.Copy($null, $wsQ3)
$wsQ3Copy = $wb.Worksheets.Item(3)

# 2) Turn the original "2022-Q3" sheet into "2022-Q4" with the new fund data
#    (rename the original out of the way first so the duplicate can reclaim
#    the "2022-Q3" name).
$wsQ3.Name = "2022-Q4"
$wsQ3Copy.Name = "2022-Q3"

# Columns B:G hold text values (fund code, name, and the numbers exactly as
# published), only H ("仓位排名") is a real number - force B:G to text first
# so values like "011924" / "90.30" keep their exact text form instead of
# being auto-coerced to numbers, then drop the helper number format again.
$textRange = $wsQ3.Range("B2:G3")
$textRange.NumberFormat = "@"

$wsQ3.Range("B2").Value = "011924"
$wsQ3.Range("C2").Value = "嘉实港股互联网产业核心资产混合A"
$wsQ3.Range("D2").Value = "1.79"
$wsQ3.Range("E2").Value = "90.30"
$wsQ3.Range("F2").Value = "4.64"
$wsQ3.Range("G2").Value = "0.0831"
$wsQ3.Range("H2").Value = 10

$wsQ3.Range("B3").Value = "011925"
$wsQ3.Range("C3").Value = "嘉实港股互联网产业核心资产混合C"
$wsQ3.Range("D3").Value = "1.14"
$wsQ3.Range("E3").Value = "90.30"
$wsQ3.Range("F3").Value = "4.64"
$wsQ3.Range("G3").Value = "0.0529"
$wsQ3.Range("H3").Value = 10

$textRange.ClearFormats()

# Header row / first column on the new "2022-Q4" sheet should carry the same
# style used on the "总计" sheet (style index 2) rather than the style that
# used to live on the old "2022-Q3" sheet (style index 1).
$wsTotal.Range("B1").Copy()
$wsQ3.Range("B1:H1").PasteSpecial(-4122)
$wsTotal.Range("A2").Copy()
$wsQ3.Range("A2:A3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 3) Update the "总计" (total) sheet: row 2 becomes 2022-Q4 and a new row 3
#    is inserted with the (unchanged) 2022-Q3 figures.
$wsTotal.Range("B2").Value = "2022-Q4"
$wsTotal.Range("D2").Value = 0.14

$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2022-Q3"
$wsTotal.Range("C3").Value = 2
$wsTotal.Range("D3").Value = 0.09

$wsTotal.Range("A2").Copy()
$wsTotal.Range("A3").PasteSpecial(-4122)

$excel.CutCopyMode = $false
